# Add new Q&A rows for three new "24 TOPICS" mini-topics (LEARNING ENGLISH,
# BIRTHDAY, MUSIC) to Sheet1, reusing the existing row-79/80 pattern and
# extending it through row 93.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column C / column D text content for rows 79-93
# ---------------------------------------------------------------------
$topics = @(
    "LEARNING ENGLISH", "LEARNING ENGLISH", "LEARNING ENGLISH", "LEARNING ENGLISH", "LEARNING ENGLISH",
    "BIRTHDAY", "BIRTHDAY", "BIRTHDAY", "BIRTHDAY", "BIRTHDAY",
    "MUSIC", "MUSIC", "MUSIC", "MUSIC", "MUSIC"
)

$questions = @(
    "1. Do you like studying English?",
    "2. How often do you study English?",
    "3. How do you study English?",
    "4. Why do you study English?",
    "5. What do you find most difficult when learning English?",
    "1. When is your birthday?",
    "2. What do you often do on your birthday?",
    "3. What gifts do you want to receive on your birthday?",
    "4. Where do you want to have your birthday party?",
    "5. What do you often have for your birthday party?",
    "1. What kind of music do you like?",
    "2. Who is your favorite singer?",
    "3. Why do you like listening to music?",
    "4. When do you listen to music?",
    "5. Where do you listen to music?"
)

$startRow = 79
for ($i = 0; $i -lt $topics.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = "24 TOPICS"
    $ws.Cells.Item($row, 2).Value = "Trả lời"
    $ws.Cells.Item($row, 3).Value = $topics[$i]
    $ws.Cells.Item($row, 4).Value = $questions[$i]
}

# ---------------------------------------------------------------------
# 2. Formatting: column C gets the plain "type" style used elsewhere in
#    the sheet (e.g. C2); column D gets a plain (non-bold, default
#    colour) Arial 10pt style instead of the old bold/coloured look.
# ---------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("C79:C93").PasteSpecial(-4122)

$helper = $ws.Range("ZZ1")
$helper.Value = "x"
$helper.Font.Name = "Arial"
$helper.Font.Size = 10
$helper.Font.Bold = $false
$helper.Font.ThemeColor = 1
$helper.Copy()
$ws.Range("D79:D93").PasteSpecial(-4122)
$helper.Clear()

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. View state: active cell / selection moves to D78 (single cell).
# ---------------------------------------------------------------------
[void]$ws.Range("D78").Select()
